$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) cells are forced to Text format before assignment so that
# numeric-looking strings (e.g. "0.440", "1.00") keep their exact textual
# representation (trailing zeros, thousand-dot separators) instead of being
# auto-converted to a Number by Excel.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '54.101.74'
$ws.Range('E2').Value = '  -7.53%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.856.76'
$ws.Range('E3').Value = '  -10.51%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '468.28'
$ws.Range('E5').Value = '  -12.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '125.45'
$ws.Range('E6').Value = '  -6.91%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.854.55'
$ws.Range('E8').Value = '  -10.53%  '
$ws.Range('E9').Value = '  -12.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.59'
$ws.Range('E10').Value = '  -11.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0954'
$ws.Range('E11').Value = '  -15.58%  '
$ws.Range('E12').Value = '  -16.85%  '
$ws.Range('E13').Value = '  -4.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.347.95'
$ws.Range('E14').Value = '  -10.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.11'
$ws.Range('E15').Value = '  -10.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '54.096.03'
$ws.Range('E16').Value = '  -7.68%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.858.06'
$ws.Range('E17').Value = '  -10.54%  '
$ws.Range('E18').Value = '  -14.84%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.32'
$ws.Range('E19').Value = '  -9.70%  '
$ws.Range('E20').Value = '  -13.79%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.99'
$ws.Range('E21').Value = '  -13.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '292.28'
$ws.Range('E22').Value = '  -18.50%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.440'
$ws.Range('E24').Value = '  -14.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '58.40'
$ws.Range('E25').Value = '  -16.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.996'
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.151'
$ws.Range('E27').Value = '  -10.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  -0.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0803'
$ws.Range('E29').Value = '  -15.44%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.13'
$ws.Range('E30').Value = '  -12.84%  '
$ws.Range('E31').Value = '  -6.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.16'
$ws.Range('E32').Value = '  -12.18%  '
$ws.Range('E33').Value = '  -16.04%  '
$ws.Range('E34').Value = '  -13.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.18'
$ws.Range('E35').Value = '  -15.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '134.11'
$ws.Range('E36').Value = '  -16.63%  '
$ws.Range('E37').Value = '  -14.87%  '
$ws.Range('E38').Value = '  -14.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '23.02'
$ws.Range('E39').Value = '  -10.68%  '
$ws.Range('E40').Value = '  -12.70%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.878.76'
$ws.Range('E41').Value = '  -10.62%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '34.69'
$ws.Range('E43').Value = '  -14.87%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.942'
$ws.Range('E44').Value = '  -14.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.598'
$ws.Range('E45').Value = '  -15.56%  '
$ws.Range('E46').Value = '  -11.90%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.36'
$ws.Range('E47').Value = '  -16.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.032.17'
$ws.Range('E48').Value = '  -10.88%  '
$ws.Range('E49').Value = '  -14.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.78'
$ws.Range('E50').Value = '  -12.98%  '
$ws.Range('E51').Value = '  -11.45%  '
